$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - add new header/label cell E1 (text "0.2", same style as D1)
$ws.Range("E1").Value = "'0.2"
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats

# Row 2
$ws.Range("A2").Value = 3.2
$ws.Range("B2").Value = 8
$ws.Range("C2").Value = 2.72727272727273
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# Row 3
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 2.22222222222222
$ws.Range("D3").Value = -0.126262626262626
$ws.Range("E3").Value = 0

# Row 4 (new row)
$ws.Range("A4").Value = -9
$ws.Range("B4").Value = -3
$ws.Range("C4").Value = 1.07142857142857
$ws.Range("D4").Value = 0.0943273484257091
$ws.Range("E4").Value = -0.0220589974688335
